# Update Name of Algo
# Update the slightly re-computed imputed values in column A (algorithm re-run results)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = -22.14960000000001
    14  = -21.76409999999999
    16  = -22.21550000000002
    21  = -20.14599999999998
    23  = -20.09859999999998
    25  = -21.73719999999998
    26  = -21.12519999999997
    29  = -20.93929999999997
    40  = -20.1001
    53  = -21.8213
    57  = -22.63960000000002
    59  = -22.6365
    65  = -21.78239999999998
    69  = -21.62069999999999
    79  = -20.12810000000001
    83  = -21.791
    91  = -20.55999999999998
    93  = -21.39450000000001
    100 = -22.17780000000001
    103 = -21.82289999999999
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
